$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stamp the workbook's existing bold header style (currently on A1) onto the
# cells that need it in the new layout, before any of the old content is
# touched. PasteSpecial(Formats) reuses the existing style record instead of
# minting a new one.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)
$ws.Range("C1:E2").PasteSpecial(-4122)
$ws.Range("C2:E2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Remove all of the old header/data content - the sheet is being
# restructured, so nothing from the old A1:G2 block survives as-is.
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").Clear()
$ws.Range("F1:G1").Clear()
$ws.Range("A2:B2").ClearContents()
$ws.Range("F2:G2").ClearContents()

# ---------------------------------------------------------------------------
# Title row ("Dealer - 09 - 2018") merged across C1:E2.
# ---------------------------------------------------------------------------
$ws.Range("C1:E2").Merge()
$ws.Range("C1").Value = "Dealer - 09 - 2018"

# ---------------------------------------------------------------------------
# Header row (row 3)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Id"
$ws.Range("B3").Value = "Date"
$ws.Range("C3").Value = "Dealer"
$ws.Range("D3").Value = "Cylinder Size"
$ws.Range("E3").Value = "Quantity"
$ws.Range("F3").Value = "Rate"
$ws.Range("G3").Value = "Net Amount"
$ws.Range("H3").Value = "Amount Paid"
$ws.Range("I3").Value = "Amount Due"

# ---------------------------------------------------------------------------
# Data row (row 4). Several numeric-looking values must be stored as TEXT
# (to match the source file), so the cells are pre-formatted as Text,
# written, then the formatting is cleared again (the Text cell type
# survives the ClearFormats call; only the number-format styling is
# removed).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 1

$ws.Range("B4:F4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"

$ws.Range("B4").Value = "2018-09-02"
$ws.Range("C4").Value = "abhishek chu"
$ws.Range("D4").Value = "12"
$ws.Range("E4").Value = "100"
$ws.Range("F4").Value = "786"
$ws.Range("G4").Value = 78600
$ws.Range("H4").Value = "7869"
$ws.Range("I4").Value = 70731

$ws.Range("B4:F4").ClearFormats()
$ws.Range("H4").ClearFormats()

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.142857142857142
$ws.Columns.Item(3).ColumnWidth = 12.142857142857142
$ws.Columns.Item(4).ColumnWidth = 16.142857142857142
$ws.Columns.Item(5).ColumnWidth = 14.142857142857142
$ws.Columns.Item(6).ColumnWidth = 9.142857142857142
$ws.Columns.Item(7).ColumnWidth = 22.142857142857142
$ws.Columns.Item(8).ColumnWidth = 22.142857142857142
$ws.Columns.Item(9).ColumnWidth = 22.142857142857142
